$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$changes = @(
    @{Cell="D2"; Value="41.736.29"}
    @{Cell="E2"; Value="  -0.70%  "}
    @{Cell="D3"; Value="2.203.39"}
    @{Cell="E3"; Value="  -2.93%  "}
    @{Cell="E4"; Value="  +0.02%  "}
    @{Cell="D5"; Value="229.84"}
    @{Cell="E5"; Value="  -2.11%  "}
    @{Cell="D6"; Value="0.616"}
    @{Cell="E6"; Value="  -4.38%  "}
    @{Cell="D7"; Value="60.16"}
    @{Cell="E7"; Value="  -5.36%  "}
    @{Cell="E8"; Value="  +0.00%  "}
    @{Cell="E9"; Value="  -2.57%  "}
    @{Cell="D10"; Value="56.98"}
    @{Cell="E10"; Value="  -5.37%  "}
    @{Cell="D11"; Value="0.0885"}
    @{Cell="E11"; Value="  -1.53%  "}
    @{Cell="E12"; Value="  -2.05%  "}
    @{Cell="D13"; Value="2.532.20"}
    @{Cell="E13"; Value="  -2.92%  "}
    @{Cell="E14"; Value="  -4.78%  "}
    @{Cell="D15"; Value="22.14"}
    @{Cell="E15"; Value="  -3.54%  "}
    @{Cell="B16"; Value="Polygon"}
    @{Cell="C16"; Value="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"}
    @{Cell="D16"; Value="0.793"}
    @{Cell="E16"; Value="  -3.74%  "}
    @{Cell="B17"; Value="Polkadot"}
    @{Cell="C17"; Value="https://coinranking.com/coin/25W7FG7om+polkadot-dot"}
    @{Cell="D17"; Value="5.56"}
    @{Cell="E17"; Value="  -2.86%  "}
    @{Cell="D18"; Value="2.215.54"}
    @{Cell="E18"; Value="  -2.35%  "}
    @{Cell="D19"; Value="41.614.15"}
    @{Cell="E19"; Value="  -0.68%  "}
    @{Cell="B20"; Value="Litecoin"}
    @{Cell="C20"; Value="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"}
    @{Cell="D20"; Value="72.03"}
    @{Cell="B21"; Value="ShibaInu"}
    @{Cell="C21"; Value="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"}
    @{Cell="D21"; Value="0.0₃0900"}
    @{Cell="E21"; Value="  -3.57%  "}
    @{Cell="E22"; Value="  -2.24%  "}
    @{Cell="D23"; Value="241.35"}
    @{Cell="E23"; Value="  -4.53%  "}
    @{Cell="E24"; Value="  -0.20%  "}
    @{Cell="D25"; Value="2.34"}
    @{Cell="E25"; Value="  -3.93%  "}
    @{Cell="E26"; Value="  -3.33%  "}
    @{Cell="D27"; Value="9.59"}
    @{Cell="E27"; Value="  -2.77%  "}
    @{Cell="D28"; Value="168.32"}
    @{Cell="E28"; Value="  -1.67%  "}
    @{Cell="E29"; Value="  -6.82%  "}
    @{Cell="E30"; Value="  -0.36%  "}
    @{Cell="D31"; Value="19.73"}
    @{Cell="E31"; Value="  -4.05%  "}
    @{Cell="D32"; Value="2.60"}
    @{Cell="E32"; Value="  -8.32%  "}
    @{Cell="D34"; Value="4.99"}
    @{Cell="E34"; Value="  -2.38%  "}
    @{Cell="E35"; Value="  -4.38%  "}
    @{Cell="D36"; Value="0.0643"}
    @{Cell="E36"; Value="  +0.45%  "}
    @{Cell="E37"; Value="  -7.43%  "}
    @{Cell="D38"; Value="6.27"}
    @{Cell="E38"; Value="  -8.38%  "}
    @{Cell="D39"; Value="2.33"}
    @{Cell="E39"; Value="  -5.36%  "}
    @{Cell="B40"; Value="BinanceUSD"}
    @{Cell="C40"; Value="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"}
    @{Cell="D40"; Value="1.00"}
    @{Cell="E40"; Value="  +0.03%  "}
    @{Cell="B41"; Value="TerraClassic"}
    @{Cell="C41"; Value="https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"}
    @{Cell="D41"; Value="0.000237"}
    @{Cell="E41"; Value="  -9.44%  "}
    @{Cell="D42"; Value="0.0240"}
    @{Cell="E42"; Value="  -1.10%  "}
    @{Cell="D43"; Value="8.63"}
    @{Cell="E43"; Value="  -0.72%  "}
    @{Cell="E44"; Value="  -14.14%  "}
    @{Cell="D45"; Value="0.0954"}
    @{Cell="E45"; Value="  -3.28%  "}
    @{Cell="D47"; Value="96.61"}
    @{Cell="E47"; Value="  -5.50%  "}
    @{Cell="D48"; Value="1.461.18"}
    @{Cell="E48"; Value="  -3.05%  "}
    @{Cell="E49"; Value="  -1.80%  "}
    @{Cell="D50"; Value="16.17"}
    @{Cell="E50"; Value="  -9.07%  "}
    @{Cell="E51"; Value="  -5.67%  "}
)

foreach ($chg in $changes) {
    $col = $chg.Cell.Substring(0,1)
    if ($col -eq "D") {
        Set-TextValue $chg.Cell $chg.Value
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}
